# Update "想去人数" (interest count) values that changed when the
# gh-pages data was regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 8943
$wsExhibit.Range("F15").Value = 18
$wsExhibit.Range("F21").Value = 1836
$wsExhibit.Range("F28").Value = 1028
$wsExhibit.Range("F31").Value = 547
$wsExhibit.Range("F34").Value = 537
$wsExhibit.Range("F35").Value = 2276
$wsExhibit.Range("F41").Value = 278
$wsExhibit.Range("F47").Value = 9

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 8943
$wsAll.Range("F19").Value = 18
$wsAll.Range("F24").Value = 1836
$wsAll.Range("F30").Value = 1028
$wsAll.Range("F32").Value = 547
$wsAll.Range("F34").Value = 537
$wsAll.Range("F35").Value = 2276
$wsAll.Range("F39").Value = 278
